$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain plain text (matches source inlineStr type),
# since Excel would otherwise auto-convert plain decimal-looking strings to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.706.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.913.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.29%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.909.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.88%  "
$ws.Range("E10").Value = "  -4.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.14%  "
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.418.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.759.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.911.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.87%  "
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.044.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.183"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.90%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.450"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0859"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.09%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("E38").Value = "  -10.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.341.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.04%  "
$ws.Range("E42").Value = "  -7.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.642"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.65%  "
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("E49").Value = "  -7.11%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0921"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.71%  "
